{"js": "// Apply the table/date updates described by the diff.\n// Each (oldText -> newText) pair is unique and unambiguous, so a plain\n// body.search()+insertText(replace) pass for every pair is sufficient.\nconst replacements = [\n  [\"2025-04-26 Saturday\", \"2025-04-27 Sunday\"],\n  [\"661\u00f73=220, 1\", \"181\u00f78=22, 5\"],\n  [\"143\u00f78=17, 7\", \"518\u00f77=74, 0\"],\n  [\"774\u00f72=387, 0\", \"196\u00f72=98, 0\"],\n  [\"444\u00f74=111, 0\", \"985\u00f78=123, 1\"],\n  [\"310\u00f78=38, 6\", \"557\u00f76=92, 5\"],\n  [\"148\u00f77=21, 1\", \"112\u00f79=12, 4\"],\n  [\"579\u00f72=289, 1\", \"126\u00f73=42, 0\"],\n  [\"799\u00f75=159, 4\", \"706\u00f77=100, 6\"],\n  [\"294\u00f72=147, 0\", \"653\u00f79=72, 5\"],\n  [\"965\u00f76=160, 5\", \"884\u00f79=98, 2\"],\n  [\"993\u00f76=165, 3\", \"174\u00f76=29, 0\"],\n  [\"579\u00f74=144, 3\", \"837\u00f79=93, 0\"],\n  [\"142\u00f79=15, 7\", \"985\u00f74=246, 1\"],\n  [\"707\u00f74=176, 3\", \"111\u00f73=37, 0\"],\n  [\"707\u00f77=101, 0\", \"292\u00f74=73, 0\"],\n  [\"300\u00f76=50, 0\", \"250\u00f74=62, 2\"],\n  [\"121\u00f77=17, 2\", \"724\u00f74=181, 0\"],\n  [\"514\u00f74=128, 2\", \"505\u00f77=72, 1\"],\n  [\"953\u00f79=105, 8\", \"174\u00f72=87, 0\"],\n  [\"690\u00f73=230, 0\", \"782\u00f78=97, 6\"],\n  [\"848\u00f74=212, 0\", \"628\u00f78=78, 4\"],\n  [\"944\u00f74=236, 0\", \"426\u00f73=142, 0\"],\n  [\"189\u00f75=37, 4\", \"532\u00f72=266, 0\"],\n  [\"623\u00f72=311, 1\", \"841\u00f72=420, 1\"],\n  [\"831\u00f73=277, 0\", \"358\u00f74=89, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the table/date updates described by the diff.\n# Each (oldText -> newText) pair is unique and unambiguous, so a plain\n# Find/Replace pass for every pair is sufficient.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-04-26 Saturday\", \"2025-04-27 Sunday\"),\n  @(\"661\u00f73=220, 1\", \"181\u00f78=22, 5\"),\n  @(\"143\u00f78=17, 7\", \"518\u00f77=74, 0\"),\n  @(\"774\u00f72=387, 0\", \"196\u00f72=98, 0\"),\n  @(\"444\u00f74=111, 0\", \"985\u00f78=123, 1\"),\n  @(\"310\u00f78=38, 6\", \"557\u00f76=92, 5\"),\n  @(\"148\u00f77=21, 1\", \"112\u00f79=12, 4\"),\n  @(\"579\u00f72=289, 1\", \"126\u00f73=42, 0\"),\n  @(\"799\u00f75=159, 4\", \"706\u00f77=100, 6\"),\n  @(\"294\u00f72=147, 0\", \"653\u00f79=72, 5\"),\n  @(\"965\u00f76=160, 5\", \"884\u00f79=98, 2\"),\n  @(\"993\u00f76=165, 3\", \"174\u00f76=29, 0\"),\n  @(\"579\u00f74=144, 3\", \"837\u00f79=93, 0\"),\n  @(\"142\u00f79=15, 7\", \"985\u00f74=246, 1\"),\n  @(\"707\u00f74=176, 3\", \"111\u00f73=37, 0\"),\n  @(\"707\u00f77=101, 0\", \"292\u00f74=73, 0\"),\n  @(\"300\u00f76=50, 0\", \"250\u00f74=62, 2\"),\n  @(\"121\u00f77=17, 2\", \"724\u00f74=181, 0\"),\n  @(\"514\u00f74=128, 2\", \"505\u00f77=72, 1\"),\n  @(\"953\u00f79=105, 8\", \"174\u00f72=87, 0\"),\n  @(\"690\u00f73=230, 0\", \"782\u00f78=97, 6\"),\n  @(\"848\u00f74=212, 0\", \"628\u00f78=78, 4\"),\n  @(\"944\u00f74=236, 0\", \"426\u00f73=142, 0\"),\n  @(\"189\u00f75=37, 4\", \"532\u00f72=266, 0\"),\n  @(\"623\u00f72=311, 1\", \"841\u00f72=420, 1\"),\n  @(\"831\u00f73=277, 0\", \"358\u00f74=89, 2\"),\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair[1]\n  $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n  if (-not $found) {\n    throw \"Text not found: $($pair[0])\"\n  }\n}\n"}
